$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated "Total" (B) and "Community" (D) values for months 1-12 (rows 2-13)
$ws.Range("B2").Value = 2993.3023958
$ws.Range("D2").Value = 202.2209602333333

$ws.Range("B3").Value = 2815.513484116667
$ws.Range("D3").Value = 193.8054265166667

$ws.Range("B4").Value = 3002.775658366667
$ws.Range("D4").Value = 200.7630757

$ws.Range("B5").Value = 2896.734895566667
$ws.Range("D5").Value = 202.4454222333333

$ws.Range("B6").Value = 2986.844631983334
$ws.Range("D6").Value = 205.8283510833333

$ws.Range("B7").Value = 2939.527072450001
$ws.Range("D7").Value = 187.4918758

$ws.Range("B8").Value = 2981.484152850001
$ws.Range("D8").Value = 195.3420448

$ws.Range("B9").Value = 3008.645255083334
$ws.Range("D9").Value = 203.6077692666667

$ws.Range("B10").Value = 2929.975138683334
$ws.Range("D10").Value = 195.5894932833333

$ws.Range("B11").Value = 3006.062622700001
$ws.Range("D11").Value = 194.65602

$ws.Range("B12").Value = 2915.079398966667
$ws.Range("D12").Value = 205.13571825

$ws.Range("B13").Value = 2931.312351233334
$ws.Range("D13").Value = 189.0492176666667
